$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Joda-time based date formulas with dateTool.format() calls that
# make use of the locale/timezone from preferences.
$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'
$ws.Range("A9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.startTime, locale, timezone)}'
$ws.Range("C9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", trip.endTime, locale, timezone)}'

# Update the selected cell in the sheet view.
$ws.Range("B2").Select()
